$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$zhcn.Range("E2").Value = "2016-03-13 22:25:43"
$zhcn.Range("H2").Value = "2016-03-13 22:26:01"
$zhcn.Range("E4").Value = "2016-03-13 22:25:43"
$zhcn.Range("H4").Value = "2016-03-13 22:26:01"

$dede.Range("E2").Value = "2016-03-13 22:25:47"
$dede.Range("H2").Value = "2016-03-13 22:26:08"
$dede.Range("E4").Value = "2016-03-13 22:25:47"
$dede.Range("H4").Value = "2016-03-13 22:26:08"
